$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update GOOGLEFINANCE fallback constants (refreshed quote values) ---
$ws.Range("C5").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("GOOGLEFINANCE(C4,$B5)"),253.47)'
$ws.Range("C6").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("GOOGLEFINANCE(C4,$B6)/1000000000"),20.414233101)'
$ws.Range("C7").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("GOOGLEFINANCE(C4,$B7)"),26.79)'

# --- Append 16 new history rows (106-121), matching the style of the last existing row (105) ---
$lastRow = 105
$firstNew = 106
$lastNew = 121
$srcFormatRange = $ws.Range("A" + $lastRow + ":C" + $lastRow)
$dstFormatRange = $ws.Range("A" + $firstNew + ":C" + $lastNew)
$srcFormatRange.Copy($dstFormatRange)

$ws.Cells.Item(106, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44739.66666666667)'
$ws.Cells.Item(106, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),248.16)'
$ws.Cells.Item(106, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),332670.0)'
$ws.Cells.Item(107, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44740.66666666667)'
$ws.Cells.Item(107, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),239.41)'
$ws.Cells.Item(107, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),320265.0)'
$ws.Cells.Item(108, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44741.66666666667)'
$ws.Cells.Item(108, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),240.14)'
$ws.Cells.Item(108, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),360665.0)'
$ws.Cells.Item(109, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44742.66666666667)'
$ws.Cells.Item(109, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),241.83)'
$ws.Cells.Item(109, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),533191.0)'
$ws.Cells.Item(110, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44743.66666666667)'
$ws.Cells.Item(110, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),244.28)'
$ws.Cells.Item(110, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),319159.0)'
$ws.Cells.Item(111, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44747.66666666667)'
$ws.Cells.Item(111, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),241.21)'
$ws.Cells.Item(111, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),366698.0)'
$ws.Cells.Item(112, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44748.66666666667)'
$ws.Cells.Item(112, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),241.77)'
$ws.Cells.Item(112, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),349285.0)'
$ws.Cells.Item(113, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44749.66666666667)'
$ws.Cells.Item(113, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),245.02)'
$ws.Cells.Item(113, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),367594.0)'
$ws.Cells.Item(114, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44750.66666666667)'
$ws.Cells.Item(114, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),242.63)'
$ws.Cells.Item(114, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),307020.0)'
$ws.Cells.Item(115, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44753.66666666667)'
$ws.Cells.Item(115, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),244.37)'
$ws.Cells.Item(115, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),439994.0)'
$ws.Cells.Item(116, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44754.66666666667)'
$ws.Cells.Item(116, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),237.75)'
$ws.Cells.Item(116, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),376703.0)'
$ws.Cells.Item(117, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44755.66666666667)'
$ws.Cells.Item(117, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),238.27)'
$ws.Cells.Item(117, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),260903.0)'
$ws.Cells.Item(118, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44756.66666666667)'
$ws.Cells.Item(118, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),235.8)'
$ws.Cells.Item(118, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),332700.0)'
$ws.Cells.Item(119, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44757.66666666667)'
$ws.Cells.Item(119, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),240.69)'
$ws.Cells.Item(119, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),338878.0)'
$ws.Cells.Item(120, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44760.66666666667)'
$ws.Cells.Item(120, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),238.77)'
$ws.Cells.Item(120, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),254069.0)'
$ws.Cells.Item(121, 1).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),44761.66666666667)'
$ws.Cells.Item(121, 2).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),247.57)'
$ws.Cells.Item(121, 3).Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),387155.0)'

Write-Host "edit complete"